$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (preserve trailing zeros / dotted formatting)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.057.17'
$ws.Range("E2").Value = '  -4.60%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.496.32'
$ws.Range("E3").Value = '  -3.22%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '536.10'
$ws.Range("E5").Value = '  -2.70%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.04'
$ws.Range("E6").Value = '  -6.82%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.995'
$ws.Range("E7").Value = '  -0.44%  '

$ws.Range("E8").Value = '  -3.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.531.88'
$ws.Range("E9").Value = '  -2.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0999'
$ws.Range("E10").Value = '  -3.90%  '

$ws.Range("E11").Value = '  -2.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.47'
$ws.Range("E12").Value = '  -0.55%  '

$ws.Range("E13").Value = '  -3.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.933.92'
$ws.Range("E14").Value = '  -3.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.83'
$ws.Range("E15").Value = '  -6.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.856.13'
$ws.Range("E16").Value = '  -4.82%  '

$ws.Range("E17").Value = '  -3.43%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.510.02'
$ws.Range("E18").Value = '  -2.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.28'
$ws.Range("E19").Value = '  -2.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.29'
$ws.Range("E20").Value = '  -5.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.07'
$ws.Range("E21").Value = '  -4.55%  '

$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.76'
$ws.Range("E23").Value = '  -4.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.70'
$ws.Range("E24").Value = '  -2.88%  '

$ws.Range("E25").Value = '  -10.43%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.161'
$ws.Range("E26").Value = '  -3.43%  '

$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.610.84'
$ws.Range("E27").Value = '  -3.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.993'
$ws.Range("E28").Value = '  -0.60%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.76'
$ws.Range("E29").Value = '  -4.53%  '

$ws.Range("E30").Value = '  -6.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0778'
$ws.Range("E31").Value = '  -6.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.20'
$ws.Range("E33").Value = '  -10.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '159.54'
$ws.Range("E34").Value = '  -2.03%  '

$ws.Range("E35").Value = '  -0.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.46'
$ws.Range("E36").Value = '  +4.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.58'
$ws.Range("E37").Value = '  -2.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.36'
$ws.Range("E38").Value = '  -10.40%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.62'
$ws.Range("E39").Value = '  -8.84%  '

$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.73'
$ws.Range("E40").Value = '  -4.89%  '

$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '305.99'
$ws.Range("E41").Value = '  -6.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '36.82'
$ws.Range("E42").Value = '  -1.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.823'
$ws.Range("E43").Value = '  -8.44%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.67'
$ws.Range("E44").Value = '  -6.25%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.994'
$ws.Range("E45").Value = '  -0.38%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.598'
$ws.Range("E46").Value = '  -1.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.76'
$ws.Range("E47").Value = '  -1.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.08'
$ws.Range("E48").Value = '  +2.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0931'
$ws.Range("E49").Value = '  -3.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.66'
$ws.Range("E50").Value = '  -4.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0518'
$ws.Range("E51").Value = '  -4.94%  '
